$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ocYZK178"
$ws.Range("B2").Value = 23071818
$ws.Range("C2").Value = "kowgmyl87"
$ws.Range("D2").Value = "ZU%3yq6$"
$ws.Range("F2").Value = "JDSAlTPD"
$ws.Range("G2").Value = "eHch"
